$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.671
$ws.Range("E2").Value = 0.599
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 0.5590936925903246

$ws.Range("D3").Value = 0.699
$ws.Range("E3").Value = 0.621
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 0.5600122473974281

$ws.Range("D4").Value = 0.717
$ws.Range("E4").Value = 0.627
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 0.629210042865891

$ws.Range("D5").Value = 0.732
$ws.Range("E5").Value = 0.629
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 0.5731781996325781

$ws.Range("D6").Value = 0.68
$ws.Range("E6").Value = 0.605
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 0.5364360073484384

$ws.Range("D7").Value = 0.71
$ws.Range("E7").Value = 0.621
$ws.Range("G7").Value = 0.5290875688916106

$ws.Range("D8").Value = 0.722
$ws.Range("E8").Value = 0.629
$ws.Range("F8").Value = 9
$ws.Range("G8").Value = 0.6175750153092467

$ws.Range("D9").Value = 0.744
$ws.Range("E9").Value = 0.633
$ws.Range("F9").Value = 16
$ws.Range("G9").Value = 0.5930802204531537

$ws.Range("D10").Value = 0.694
$ws.Range("E10").Value = 0.616
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 0.5263319044703001

$ws.Range("D11").Value = 0.712
$ws.Range("E11").Value = 0.623
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 0.5835884874464177

$ws.Range("D12").Value = 0.727
$ws.Range("E12").Value = 0.633
$ws.Range("G12").Value = 0.6313533374157991

$ws.Range("D13").Value = 0.75
$ws.Range("E13").Value = 0.635
$ws.Range("G13").Value = 0.5783833435394978

$ws.Range("D14").Value = 0.6929999999999999
$ws.Range("E14").Value = 0.614
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 0.5554194733619106

$ws.Range("D15").Value = 0.712
$ws.Range("G15").Value = 0.5315370483772198

$ws.Range("D16").Value = 0.736
$ws.Range("E16").Value = 0.634
$ws.Range("G16").Value = 0.5627679118187385

$ws.Range("D17").Value = 0.754
$ws.Range("E17").Value = 0.638
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 0.5774647887323944

$ws.Range("D18").Value = 0.701
$ws.Range("E18").Value = 0.623
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 0.5759338640538886

$ws.Range("E19").Value = 0.627
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 0.5202082057562768

$ws.Range("D20").Value = 0.743
$ws.Range("E20").Value = 0.636
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 0.6062461726883037

$ws.Range("D21").Value = 0.762
$ws.Range("E21").Value = 0.642
$ws.Range("F21").Value = 22
$ws.Range("G21").Value = 0.580526638089406

$ws.Range("D22").Value = 0.717
$ws.Range("E22").Value = 0.627
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 0.5569503980404165

$ws.Range("D23").Value = 0.735
$ws.Range("E23").Value = 0.633
$ws.Range("G23").Value = 0.5535823637477036

$ws.Range("D24").Value = 0.755
$ws.Range("E24").Value = 0.643
$ws.Range("F24").Value = 17
$ws.Range("G24").Value = 0.5970606246172688

$ws.Range("D25").Value = 0.77
$ws.Range("E25").Value = 0.644
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 0.5725658297611758

$ws.Range("D26").Value = 0.722
$ws.Range("E26").Value = 0.63
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 0.5731781996325781

$ws.Range("D27").Value = 0.742
$ws.Range("E27").Value = 0.639
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 0.5496019595835885

$ws.Range("D28").Value = 0.76
$ws.Range("E28").Value = 0.643
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 0.5946111451316595

$ws.Range("D29").Value = 0.774
$ws.Range("E29").Value = 0.645
$ws.Range("G29").Value = 0.575015309246785

$ws.Range("D30").Value = 0.727
$ws.Range("G30").Value = 0.5759338640538886

$ws.Range("D31").Value = 0.743
$ws.Range("E31").Value = 0.639
$ws.Range("F31").Value = 14
$ws.Range("G31").Value = 0.5673606858542559

$ws.Range("D32").Value = 0.764
$ws.Range("E32").Value = 0.646
$ws.Range("F32").Value = 19
$ws.Range("G32").Value = 0.5780771586037967

$ws.Range("D33").Value = 0.781
$ws.Range("E33").Value = 0.649
$ws.Range("F33").Value = 29
$ws.Range("G33").Value = 0.5688916105327618

$ws.Range("D34").Value = 0.73
$ws.Range("E34").Value = 0.635
$ws.Range("G34").Value = 0.5740967544396816

$ws.Range("D35").Value = 0.751
$ws.Range("E35").Value = 0.643
$ws.Range("F35").Value = 14
$ws.Range("G35").Value = 0.5673606858542559

$ws.Range("D36").Value = 0.767
$ws.Range("E36").Value = 0.647
$ws.Range("F36").Value = 22
$ws.Range("G36").Value = 0.5832823025107164

$ws.Range("D37").Value = 0.787
$ws.Range("E37").Value = 0.653
$ws.Range("F37").Value = 30
$ws.Range("G37").Value = 0.5673606858542559

$ws.Range("D38").Value = 0.736
$ws.Range("E38").Value = 0.635
$ws.Range("G38").Value = 0.5808328230251072

$ws.Range("D39").Value = 0.756
$ws.Range("E39").Value = 0.644
$ws.Range("G39").Value = 0.5646050214329456

$ws.Range("D40").Value = 0.775
$ws.Range("E40").Value = 0.649
$ws.Range("F40").Value = 23
$ws.Range("G40").Value = 0.5802204531537049

$ws.Range("D41").Value = 0.795
$ws.Range("E41").Value = 0.653
$ws.Range("F41").Value = 31
$ws.Range("G41").Value = 0.5621555419473362

$ws.Range("D42").Value = 0.743
$ws.Range("E42").Value = 0.639
$ws.Range("G42").Value = 0.5731781996325781

$ws.Range("D43").Value = 0.765
$ws.Range("E43").Value = 0.648
$ws.Range("G43").Value = 0.5701163502755664

$ws.Range("D44").Value = 0.783
$ws.Range("E44").Value = 0.654
$ws.Range("F44").Value = 27
$ws.Range("G44").Value = 0.5820575627679119

$ws.Range("D45").Value = 0.8080000000000001
$ws.Range("E45").Value = 0.656
$ws.Range("F45").Value = 31
$ws.Range("G45").Value = 0.546846295162278

$ws.Range("D46").Value = 0.752
$ws.Range("E46").Value = 0.644
$ws.Range("G46").Value = 0.5725658297611758

$ws.Range("D47").Value = 0.778
$ws.Range("E47").Value = 0.649
$ws.Range("F47").Value = 17
$ws.Range("G47").Value = 0.5676668707899571

$ws.Range("D48").Value = 0.797
$ws.Range("E48").Value = 0.657
$ws.Range("F48").Value = 28
$ws.Range("G48").Value = 0.5639926515615432

$ws.Range("D49").Value = 0.819
$ws.Range("E49").Value = 0.659
$ws.Range("G49").Value = 0.533680342927128

